$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: CheckID 8 - Sleeping Query with Open Transactions
$ws.Range("A12").Value = 8
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = "Query Problems"
$ws.Range("D12").Value = "Sleeping Query with Open Transactions"
$ws.Range("E12").Value = "http://BrentOzar.com/go/sleeping"
$ws.Hyperlinks.Add($ws.Range("E12"), "http://BrentOzar.com/go/sleeping") | Out-Null

# Row 13: CheckID 9 - Query Rolling Back
$ws.Range("A13").Value = 9
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = "Query Problems"
$ws.Range("D13").Value = "Query Rolling Back"
$ws.Range("E13").Value = "http://BrentOzar.com/go/rollback"
$ws.Hyperlinks.Add($ws.Range("E13"), "http://BrentOzar.com/go/rollback") | Out-Null

$ws.Range("E13").Select()
